$d = $word.ActiveDocument

# Locate the paragraph that reads "Core: TBX-Core" (the Core module line).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Core: TBX-Core*") {
        $target = $cand
        break
    }
}

$r = $target.Range
$full = $r.Text
$idx = $full.IndexOf("TBX-")
$delStart = $r.Start + $idx
$delEnd = $delStart + 4

# Remove the stray "TBX-" text so the line reads "Core: Core".
$d.Range($delStart, $delEnd).Delete()

# Re-locate the paragraph after the edit and figure out where the
# (now adjacent) "Core" run begins, i.e. right after ": ".
$p2 = $target.Range
$full2 = $p2.Text
$bmPos = $p2.Start + $full2.IndexOf(": ") + 2

# Move the "_GoBack" bookmark here (Word automatically relocates a
# bookmark of the same name rather than creating a duplicate).
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
